# Append the new game row (row 2) to the "Jogos do Dia" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns -----------------------------------------------------
$ws.Range("A2").Value = "FIFA World Cup Qualifiers - Asia"

# B2 holds a date-like string ("2025-11-13") that must stay plain text
# (it must NOT be auto-converted into an Excel date serial number), so we
# force a text number format before assigning it, then restore the cell
# to the default "Normal" style so no stray style index is left behind.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-11-13"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "13:00:00"
$ws.Range("D2").Value = "UAE"
$ws.Range("E2").Value = "Iraq"

# --- Numeric odds columns ----------------------------------------------
$ws.Range("F2").Value = 2.1
$ws.Range("G2").Value = 2.2
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 4.6
$ws.Range("J2").Value = 3.3
$ws.Range("K2").Value = 3.4
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 2.42
$ws.Range("O2").Value = 1.58
$ws.Range("P2").Value = 1.49
$ws.Range("Q2").Value = 2.68
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 6
$ws.Range("T2").Value = 2.22
$ws.Range("U2").Value = 1.68
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 1.67
$ws.Range("X2").Value = 8.6
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 30
$ws.Range("AA2").Value = 120
$ws.Range("AB2").Value = 6.6
$ws.Range("AC2").Value = 7.8
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 95
$ws.Range("AF2").Value = 11.5
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 150
$ws.Range("AJ2").Value = 30
$ws.Range("AK2").Value = 34
$ws.Range("AL2").Value = 75
$ws.Range("AM2").Value = 250
$ws.Range("AN2").Value = 34
$ws.Range("AO2").Value = 150
